$d = $word.ActiveDocument

function Replace-InParagraph($index, $find, $replace) {
    $p = $d.Paragraphs.Item($index).Range
    $ok = $p.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host ("FAILED to replace in paragraph " + $index + ": [" + $find + "]")
    }
}

# Paragraph 1: hyperlink "English" -> "Englisch", and language list
Replace-InParagraph 1 "English" "Englisch"
Replace-InParagraph 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugiesisch / Französisch / Thailändisch / Vietnamesisch / Spanisch"

# Paragraph 3: "English" -> "Englisch"
Replace-InParagraph 3 "English" "Englisch"

# Paragraph 6: Brief description
Replace-InParagraph 6 "An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io" "Eine E-Mail, die an Partner im Zielland gesendet wird, die mit 'Ja' geantwortet haben, aber uns ihre Dokumente nicht gesendet haben. Es wird über customer.io gesendet"

# Paragraph 8: "Target audience" -> "Zielgruppe"
Replace-InParagraph 8 "Target audience" "Zielgruppe"

# Paragraph 9: invited partners
Replace-InParagraph 9 "Invited partners who haven’t submitted their documents" "Eingeladene Partner, die ihre Dokumente noch nicht eingereicht haben"

# Paragraph 12: Subject line (informal template)
Replace-InParagraph 12 "Subject line" "Betreffzeile"
Replace-InParagraph 12 " — have you submitted your docs?  " " — haben Sie Ihre Dokumente eingereicht?  "

# Paragraph 14: Don't forget heading (informal)
Replace-InParagraph 14 "Don’t forget to send your documents" "Vergessen Sie nicht, Ihre Dokumente zu schicken"

# Paragraph 16: "Hi " -> "Hallo "
Replace-InParagraph 16 "Hi " "Hallo "

# Paragraph 18: excited to see you at upcoming event
Replace-InParagraph 18 "We’re excited to see you at the upcoming " "Wir freuen uns, Sie bei der nächsten "
Replace-InParagraph 18 ". " " zu sehen. "

# Paragraph 19: confirm registration
Replace-InParagraph 19 "To confirm your registration, we need the following documents from you by " "Um Ihre Anmeldung zu bestätigen, benötigen wir von Ihnen bis zum "
Replace-InParagraph 19 ":" " die folgenden Dokumente:"

# Paragraph 20: insert list of documents required
Replace-InParagraph 20 "[insert list of documents required]" "[Liste der erforderlichen Dokumente einfügen]"

# Paragraph 21: please send a copy
Replace-InParagraph 21 "Please send a copy of these documents to your country manager, " "Bitte senden Sie eine Kopie dieser Dokumente an Ihren Country Manager, "
Replace-InParagraph 21 ", at " ", unter "
Replace-InParagraph 21 " or " " oder "
Replace-InParagraph 21 " (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation." " (WhatsApp), damit wir die notwendigen Vorkehrungen für Sie treffen können, einschließlich Unterkunft und Transport."

# Paragraph 22: If you have any questions contact country manager
Replace-InParagraph 22 "If you have any questions, please contact your country manager." "Wenn Sie Fragen haben, wenden Sie sich bitte an Ihren Ländermanager."

# Paragraph 23: We look forward
Replace-InParagraph 23 "We look forward to seeing you there!" "Wir freuen uns darauf, Sie dort zu sehen!"

# Paragraph 29: Subject line (formal template)
Replace-InParagraph 29 "Subject line" "Betreffzeile"
Replace-InParagraph 29 " — have you submitted your docs?  " " — haben Sie Ihre Dokumente eingereicht?  "

# Paragraph 31: Don't forget heading (formal)
Replace-InParagraph 31 "Don’t forget to send your documents" "Vergessen Sie nicht, Ihre Dokumente zu schicken"

# Paragraph 33: "Dear " -> "Sehr geehrter "
Replace-InParagraph 33 "Dear " "Sehr geehrter "

# Paragraph 35: excited to see you (formal, note trailing curly quote)
Replace-InParagraph 35 "We’re excited to see you at the upcoming " "Wir freuen uns, Sie bei der nächsten "
Replace-InParagraph 35 ". ‘" " zu sehen. ‘"

# Paragraph 36: to ensure best experience
Replace-InParagraph 36 "To ensure you have the best experience at this event, we need the following documents from you by " "Um sicherzustellen, dass Sie die beste Erfahrung bei dieser Veranstaltung machen, benötigen wir von Ihnen bis zum "
Replace-InParagraph 36 ":" " die folgenden Dokumente:"

# Paragraph 37: insert list of documents required (formal)
Replace-InParagraph 37 "[insert list of documents required]" "[Liste der erforderlichen Dokumente einfügen]"

# Paragraph 38: please reply to this email
Replace-InParagraph 38 "Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation." "Bitte antworten Sie auf diese E-Mail mit einer Kopie dieser Dokumente, damit wir die notwendigen Vorkehrungen für Sie treffen können, einschließlich Unterkunft und Transport."

# Paragraph 39: if you have any questions via live chat / WhatsApp
Replace-InParagraph 39 "If you have any questions, please contact us via " "Wenn Sie Fragen haben, wenden Sie sich bitte über "
Replace-InParagraph 39 "live chat" "Live-Chat"
Replace-InParagraph 39 " or " " oder "
Replace-InParagraph 39 ". " " an uns. "

# Paragraph 40: if you have any questions contact country manager, NAME
Replace-InParagraph 40 "If you have any questions, please contact your country manager, " "Wenn Sie Fragen haben, wenden Sie sich bitte an Ihren Country Manager, "
Replace-InParagraph 40 ", at " ", unter "
Replace-InParagraph 40 " or " " oder "

# Paragraph 41: We look forward (formal)
Replace-InParagraph 41 "We look forward to seeing you there!" "Wir freuen uns darauf, Sie dort zu sehen!"

# Comment: "choose either one" -> "Wählen Sie eines davon aus"
$c1 = $d.Comments.Item(1)
$c1.Range.Text = "Wählen Sie eines davon aus"
